# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the d718df71-... files in both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status text for the d718df71 file changes everywhere it is shown: the
# Overview summary columns (zh-cn/de-de) as well as the per-language Status
# column on each language sheet (row 3 = the d718df71 file in both tables).
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# zh-cn sheet: widen the Error Detail column (P) and set the error message
# for the d718df71 row (row 3)
# (39.1666... is chosen so the stored OOXML column width attribute comes out
# to exactly 40, matching the target workbook's "<col ... width="40" />".)
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsZhCn.Range("P3").Value = "Handback file name: eecgwotr.un3 is different with handoff file name: d718df71-e7e9-4531-9f8a-fc895e7840a5.bf8bdc8d3762afcb0b529bd355db74e9f536944c.zh-cn."

# de-de sheet: widen the Error Detail column (P) and set the error message
# for the d718df71 row (row 3)
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsDeDe.Range("P3").Value = "Handback file name: eecgwotr.un3 is different with handoff file name: d718df71-e7e9-4531-9f8a-fc895e7840a5.bf8bdc8d3762afcb0b529bd355db74e9f536944c.de-de."
